# Update cfb_weather.xlsx with Timestamp 2024-10-18T10:03:01.195107
# and refreshed weather/odds figures pulled at the new run time.

$wb = $excel.ActiveWorkbook
$fbs = $wb.Worksheets.Item("FBS")
$other = $wb.Worksheets.Item("Other")

# --- Timestamp refresh (every row on the FBS sheet shares this string) ---
$fbs.Range("AK2:AK57").Value = "2024-10-18T10:03:01.195107"

# --- FBS sheet: row 4 (Baylor @ Texas Tech) ---
$fbs.Range("Y4").Value = 55.5
$fbs.Range("Z4").Value = -115
$fbs.Range("AE4").Value = -0.03478260869565217

# --- FBS sheet: row 6 (Florida State @ Duke) ---
$fbs.Range("N6").Value = "S"
$fbs.Range("O6").Value = 53.9
$fbs.Range("P6").Value = 5
$fbs.Range("Q6").Value = "S"
$fbs.Range("U6").Value = 0.2

# --- FBS sheet: row 7 (Oregon @ Purdue) ---
$fbs.Range("O7").Value = 51.85999999999999
$fbs.Range("P7").Value = 3.8
$fbs.Range("U7").Value = -6.5

# --- FBS sheet: row 8 (Oklahoma State @ Brigham Young) ---
$fbs.Range("O8").Value = 45.13999999999999
$fbs.Range("P8").Value = 5.1
$fbs.Range("Q8").Value = "NNE"
$fbs.Range("U8").Value = -0.4

# --- FBS sheet: row 9 (Fresno State @ Nevada) ---
$fbs.Range("N9").Value = "SSW"
$fbs.Range("O9").Value = 42.68
$fbs.Range("P9").Value = 5.9
$fbs.Range("Q9").Value = "SSW"
$fbs.Range("U9").Value = 0.3

# --- FBS sheet: row 11 (East Carolina @ Army) ---
$fbs.Range("Z11").Value = -115

# --- FBS sheet: row 26 (Charlotte @ Navy) ---
$fbs.Range("Q26").Value = "SSW"

# --- FBS sheet: row 27 (Alabama @ Tennessee) ---
$fbs.Range("Q27").Value = "WSW"

# --- FBS sheet: row 39 (LSU @ Arkansas) ---
$fbs.Range("AB39").Value = 2.5
$fbs.Range("AF39").Value = 0

# --- FBS sheet: row 42 (Kansas State @ West Virginia) ---
$fbs.Range("Q42").Value = "NNE"

# --- FBS sheet: row 48 (SMU @ Stanford) ---
$fbs.Range("Z48").Value = -118

# --- FBS sheet: row 54 (Western Kentucky @ Sam Houston State) ---
$fbs.Range("Q54").Value = "WSW"

# --- FBS sheet: row 56 (Georgia State @ Marshall) ---
$fbs.Range("Q56").Value = "NNE"

# --- Other sheet: row 5 (Brown vs Princeton) ---
$other.Range("O5").Value = "SW"
$other.Range("P5").Value = "SW"
$other.Range("R5").Value = 5.7
$other.Range("S5").Value = "SW"

# --- Other sheet: row 9 (Stony Brook vs Towson) ---
$other.Range("S9").Value = "SE"

# --- Other sheet: row 21 (Wofford vs Chattanooga) ---
$other.Range("S21").Value = "WSW"

# --- Other sheet: row 39 (South Dakota vs Youngstown State) ---
$other.Range("S39").Value = "SSW"
